# Apply "added some notes for tomorrow's lab" edit:
#  1. Remove the stray _GoBack bookmark that was sitting after
#     ", back down to 24" inside the table.
#  2. Append a new trailing paragraph at the very end of the document
#     (after the existing blank paragraph, before the sectPr) containing
#     the 11/25/2019 update note, with the "Update 11/25/2019: " lead-in
#     in bold and the remainder in regular weight.
#  3. Re-anchor the _GoBack bookmark at the end of that new paragraph,
#     which is where Word leaves it after the last edit.

$d = $word.ActiveDocument

# --- 1. Drop the old _GoBack bookmark from the table cell paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Add the new trailing paragraph with the update note ---
$lastPara = $d.Paragraphs.Last
$tail = $lastPara.Range
$tail.InsertParagraphAfter()

$notePara = $d.Paragraphs.Last
$noteRange = $notePara.Range

# Make sure the new paragraph carries the same Times New Roman / 12pt
# formatting used throughout the rest of the document.
$noteRange.Font.Name = "Times New Roman"
$noteRange.Font.Size = 12

$leadIn = "Update 11/25/2019: "
$bodyText = "Backed down to BASE_SPEED = 100 and Kp (6) and Kp (24). " + `
    "Near perfect run. What we need to do is get rid of the 100ms delay " + `
    "after the car does the doughnut. Since the sensors are not lined up " + `
    "with the axle, once it flips it will be off the tape anyway, " + `
    "therefore we don" + [char]0x2019 + "t need the offset push for 100ms. " + `
    "If there is time, recommend tripling speed to 180 and testing until perfect."

# Type the whole note in one shot (keeps the Range bookkeeping simple/stable),
# then go back and bold just the "Update 11/25/2019: " lead-in.
$noteStart = $noteRange.Start
$noteRange.InsertAfter($leadIn + $bodyText)

$boldRange = $d.Range($noteStart, $noteStart + $leadIn.Length)
$boldRange.Font.Bold = 1

# --- 3. Leave the _GoBack bookmark at the end of the new note, which is
#        where the cursor was after the last edit. ---
$endRange = $notePara.Range
$endRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRange)
